$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): keep A1, change B1:D1 labels, add E1:F1 ---
# Give the new header cells (E1:F1) the same style as the existing header
# cells (bold font, border, centered) by copying D1's formatting first.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("B1").Value = "Question 1"
$ws.Range("C1").Value = "Question 2"
$ws.Range("D1").Value = "Question 3"
$ws.Range("E1").Value = "Question 4"
$ws.Range("F1").Value = "Question 5"

# --- Row 2: update existing record, add new columns E2:F2 ---
$ws.Range("B2").Value = "Test user"
$ws.Range("C2").Value = 55
$ws.Range("D2").Value = "Mitsubishi"
$ws.Range("E2").Value = "moscow"
$ws.Range("F2").Value = 500

# --- Row 3: new record ---
$ws.Range("A3").Value = 396358609
$ws.Range("B3").Value = "oleg"
$ws.Range("C3").Value = 44
$ws.Range("D3").Value = "mazda"
$ws.Range("E3").Value = "rovno"
$ws.Range("F3").Value = 200

# --- Row 4: new record ---
$ws.Range("A4").Value = 396358610
$ws.Range("B4").Value = "ggg"
$ws.Range("C4").Value = 55
$ws.Range("D4").Value = "toyota"
$ws.Range("E4").Value = "wroclaw"
$ws.Range("F4").Value = 500

# --- Row 5: new record. C5 and F5 are stored as text (not numbers) ---
$ws.Range("A5").Value = 396358608
$ws.Range("B5").Value = "test driver"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "11"
$ws.Range("D5").Value = "logan"
$ws.Range("E5").Value = "rostow"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "20"
